# Quarterly financials update for VRSN_QTR_FIN.xlsx
# Inserts two new quarterly columns (D:E) ahead of the existing data,
# shifting the historical quarters right, fills in the two new quarters
# of data, and corrects a handful of historical figures that were
# revised in this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two blank columns before column D; everything from D.. shifts
#    right by two (old D -> F, ... old K -> M).
$ws.Range("D:E").Insert(-4161) | Out-Null

# 2) Copy number formats (date / #,##0) from column F (the old column D)
#    into the two freshly inserted columns so the new quarters render the
#    same way as the rest of the table.
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the two new quarters of data.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 307500
$ws.Range("E8").Value2 = 305800
$ws.Range("D9").Value2 = 48400
$ws.Range("E9").Value2 = 48200
$ws.Range("D10").Value2 = 259100
$ws.Range("E10").Value2 = 257600
$ws.Range("D12").Value2 = 15000
$ws.Range("E12").Value2 = 13700
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = -54800
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 58600
$ws.Range("E17").Value2 = 110800
$ws.Range("D18").Value2 = 248900
$ws.Range("E18").Value2 = 195000
$ws.Range("D20").Value2 = 7700
$ws.Range("E20").Value2 = 5900
$ws.Range("D21").Value2 = 268500
$ws.Range("E21").Value2 = 213200
$ws.Range("D22").Value2 = 22600
$ws.Range("E22").Value2 = 22600
$ws.Range("D23").Value2 = 233900
$ws.Range("E23").Value2 = 178300
$ws.Range("D24").Value2 = 46000
$ws.Range("E24").Value2 = 40600
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 187900
$ws.Range("E26").Value2 = 137700
$ws.Range("D27").Value2 = 187900
$ws.Range("E27").Value2 = 137700
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = -5700
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -7700
$ws.Range("E32").Value2 = -5900
$ws.Range("D33").Value2 = 182200
$ws.Range("E33").Value2 = 137700
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 182200
$ws.Range("E35").Value2 = 137700
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 357400
$ws.Range("E41").Value2 = 231600
$ws.Range("D42").Value2 = 912300
$ws.Range("E42").Value2 = 947400
$ws.Range("D43").Value2 = 10500
$ws.Range("E43").Value2 = 12100
$ws.Range("D44").Value2 = 0
$ws.Range("E44").Value2 = 0
$ws.Range("D45").Value2 = 36900
$ws.Range("E45").Value2 = 44900
$ws.Range("D46").Value2 = 1317000
$ws.Range("E46").Value2 = 1236000
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 253900
$ws.Range("E48").Value2 = 256300
$ws.Range("D49").Value2 = 52500
$ws.Range("E49").Value2 = 52500
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 291000
$ws.Range("E52").Value2 = 339900
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 1914500
$ws.Range("E54").Value2 = 1884600
$ws.Range("D57").Value2 = 20300
$ws.Range("E57").Value2 = 18600
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("D59").Value2 = 927300
$ws.Range("E59").Value2 = 895000
$ws.Range("D60").Value2 = 947600
$ws.Range("E60").Value2 = 913600
$ws.Range("D61").Value2 = 1785000
$ws.Range("E61").Value2 = 1784400
$ws.Range("D62").Value2 = 567300
$ws.Range("E62").Value2 = 587700
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 3300000
$ws.Range("E66").Value2 = 3285700
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = -17089800
$ws.Range("E72").Value2 = -17272000
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = -1385500
$ws.Range("E76").Value2 = -1401100
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 182200
$ws.Range("E81").Value2 = 137700
$ws.Range("D83").Value2 = 11900
$ws.Range("E83").Value2 = 12300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 218500
$ws.Range("E89").Value2 = 187500
$ws.Range("D91").Value2 = -7400
$ws.Range("E91").Value2 = -10900
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = 85700
$ws.Range("E94").Value2 = -38700
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -178300
$ws.Range("E100").Value2 = -173200
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = -400
$ws.Range("D102").Value2 = 125900
$ws.Range("E102").Value2 = -24800

# 4) A few historical figures were also restated in this update (now in
#    column H after the two-column insert).
$ws.Range("H24").Value2 = 30000
$ws.Range("H26").Value2 = 112000
$ws.Range("H27").Value2 = 112000
$ws.Range("H29").Value2 = -9200

# 5) Re-fit the two new data columns so their width matches the rest of
#    the table.
$ws.Range("D:E").EntireColumn.AutoFit() | Out-Null

Write-Host ("Done. Used range: " + $ws.UsedRange.Address())
